# Generate Report for Handoff
# Adds a new file entry (3c66cc8d-ef59-4e25-a0ea-1e850a972202.md and its
# localized handoff artifacts) ahead of the existing 9a6ee149... entry on
# every worksheet, and pushes the ".localization-config" row down.

$wb = $excel.ActiveWorkbook

$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/75213c72128d93fd5147981c47b55341922ffca2"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9c779d0c3c1f94d23684cf067f1aabb560057c0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/88b4fc4a836daa6211798f12c9f1c12620a3d50c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$newFile = "3c66cc8d-ef59-4e25-a0ea-1e850a972202.md"
$oldFile = "9a6ee149-ddeb-48da-bb21-231154f7bf27.md"
$cfgFile = ".localization-config"

$newZhXlf = "3c66cc8d-ef59-4e25-a0ea-1e850a972202.e0db2fdc3b7bb579d67b7b59706e2c775a40573d.zh-cn.xlf"
$newDeXlf = "3c66cc8d-ef59-4e25-a0ea-1e850a972202.e0db2fdc3b7bb579d67b7b59706e2c775a40573d.de-de.xlf"
$oldZhXlf = "9a6ee149-ddeb-48da-bb21-231154f7bf27.9da0733e958132f05a35bbfed047050952dac720.zh-cn.xlf"
$oldDeXlf = "9a6ee149-ddeb-48da-bb21-231154f7bf27.9da0733e958132f05a35bbfed047050952dac720.de-de.xlf"

$readyStatus = "Ready for handoff"
$notLocStatus = "Not to be localized"
$includeStatus = "Include"
$ignoredStatus = "Ignored"
$epoch = "0001-01-01 00:00:00"

$newZhDate = "2016-03-09 02:52:05"
$oldZhDate = "2016-03-09 02:51:13"
$newDeDate = "2016-03-09 02:52:15"
$oldDeDate = "2016-03-09 02:51:24"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" : File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Cells.Item(2,1).Value = $newFile
$ws1.Cells.Item(2,2).Value = $readyStatus
$ws1.Cells.Item(2,3).Value = $readyStatus

$ws1.Cells.Item(3,1).Value = $oldFile
$ws1.Cells.Item(3,2).Value = $readyStatus
$ws1.Cells.Item(3,3).Value = $readyStatus

$ws1.Cells.Item(4,1).Value = $cfgFile
$ws1.Cells.Item(4,2).Value = $notLocStatus
$ws1.Cells.Item(4,3).Value = $notLocStatus

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "$repoBase/e2e/$newFile", "", "", $newFile)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$repoBase/e2e/$oldFile", "", "", $oldFile)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$repoBase/$cfgFile", "", "", $cfgFile)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Cells.Item(2,1).Value = $newFile
$ws2.Cells.Item(2,2).Value = $readyStatus
$ws2.Cells.Item(2,3).Value = $newZhXlf
$ws2.Cells.Item(2,4).NumberFormat = $dateFmt
$ws2.Cells.Item(2,4).Value = $newZhDate
$ws2.Cells.Item(2,7).NumberFormat = $dateFmt
$ws2.Cells.Item(2,7).Value = $epoch
$ws2.Cells.Item(2,8).Value = $includeStatus

$ws2.Cells.Item(3,1).Value = $oldFile
$ws2.Cells.Item(3,2).Value = $readyStatus
$ws2.Cells.Item(3,3).Value = $oldZhXlf
$ws2.Cells.Item(3,4).NumberFormat = $dateFmt
$ws2.Cells.Item(3,4).Value = $oldZhDate
$ws2.Cells.Item(3,7).NumberFormat = $dateFmt
$ws2.Cells.Item(3,7).Value = $epoch
$ws2.Cells.Item(3,8).Value = $includeStatus

$ws2.Cells.Item(4,1).Value = $cfgFile
$ws2.Cells.Item(4,2).Value = $notLocStatus
$ws2.Cells.Item(4,4).NumberFormat = $dateFmt
$ws2.Cells.Item(4,4).Value = $epoch
$ws2.Cells.Item(4,7).NumberFormat = $dateFmt
$ws2.Cells.Item(4,7).Value = $epoch
$ws2.Cells.Item(4,8).Value = $ignoredStatus

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "$repoBase/e2e/$newFile", "", "", $newFile)
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhHandoffBase/$newZhXlf", "", "", $newZhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$repoBase/e2e/$oldFile", "", "", $oldFile)
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhHandoffBase/$oldZhXlf", "", "", $oldZhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$repoBase/$cfgFile", "", "", $cfgFile)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Cells.Item(2,1).Value = $newFile
$ws3.Cells.Item(2,2).Value = $readyStatus
$ws3.Cells.Item(2,3).Value = $newDeXlf
$ws3.Cells.Item(2,4).NumberFormat = $dateFmt
$ws3.Cells.Item(2,4).Value = $newDeDate
$ws3.Cells.Item(2,7).NumberFormat = $dateFmt
$ws3.Cells.Item(2,7).Value = $epoch
$ws3.Cells.Item(2,8).Value = $includeStatus

$ws3.Cells.Item(3,1).Value = $oldFile
$ws3.Cells.Item(3,2).Value = $readyStatus
$ws3.Cells.Item(3,3).Value = $oldDeXlf
$ws3.Cells.Item(3,4).NumberFormat = $dateFmt
$ws3.Cells.Item(3,4).Value = $oldDeDate
$ws3.Cells.Item(3,7).NumberFormat = $dateFmt
$ws3.Cells.Item(3,7).Value = $epoch
$ws3.Cells.Item(3,8).Value = $includeStatus

$ws3.Cells.Item(4,1).Value = $cfgFile
$ws3.Cells.Item(4,2).Value = $notLocStatus
$ws3.Cells.Item(4,4).NumberFormat = $dateFmt
$ws3.Cells.Item(4,4).Value = $epoch
$ws3.Cells.Item(4,7).NumberFormat = $dateFmt
$ws3.Cells.Item(4,7).Value = $epoch
$ws3.Cells.Item(4,8).Value = $ignoredStatus

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "$repoBase/e2e/$newFile", "", "", $newFile)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deHandoffBase/$newDeXlf", "", "", $newDeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$repoBase/e2e/$oldFile", "", "", $oldFile)
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deHandoffBase/$oldDeXlf", "", "", $oldDeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$repoBase/$cfgFile", "", "", $cfgFile)

"Done generating handoff report rows"
